$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held some values as text (shared strings "0.5", "0.004",
# "0.3", "0.002"). After the review they become real numeric values, and a
# few of them (F2/G2/M2/N2) were corrected to different numbers entirely.
$ws.Range("A2").Value = 0.5
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 80
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 0.3
$ws.Range("I2").Value = 0.3
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 8

# Apply a 3-decimal numeric format to the whole data row (new style, based
# on the existing font-1 style used by these cells).
$ws.Range("A2:N2").NumberFormat = "0.000"

# The selection left behind after the review session.
[void]$ws.Range("N3").Select()
